$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 775
$ws.Range("J32").Value = 775
$ws.Range("L32").Value = 775
$ws.Range("N32").Value = -1427
$ws.Range("H33").Value = 496.92307
$ws.Range("I33").Value = 500.85715
$ws.Range("K33").Value = 500.85715
$ws.Range("M33").Value = -271.85715
$ws.Range("H87").Value = 25438.166
$ws.Range("J87").Value = 25438.166
$ws.Range("L87").Value = 25438.166
$ws.Range("N87").Value = -27934.166
$ws.Range("H90").Value = 25438.166
$ws.Range("J90").Value = 25438.166
$ws.Range("L90").Value = 76314.49800000001
$ws.Range("N90").Value = -88794.49800000001
$ws.Range("H113").Value = 4681.154
$ws.Range("J113").Value = 6489
$ws.Range("L113").Value = 6489
$ws.Range("N113").Value = -12997
$ws.Range("H138").Value = 2798.5671
$ws.Range("I138").Value = 1892.5
$ws.Range("J138").Value = 3732.0908
$ws.Range("K138").Value = 5677.5
$ws.Range("L138").Value = 11196.2724
$ws.Range("M138").Value = -537.5
$ws.Range("N138").Value = -21476.2724

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4903352.5
$ws.Range("I2").Value = 8929646
$ws.Range("J2").Value = 1776.9131
$ws.Range("K2").Value = 8929646
$ws.Range("L2").Value = 1776.9131
$ws.Range("M2").Value = -8929533
$ws.Range("N2").Value = -2002.9131
$ws.Range("H45").Value = 1545.2069
$ws.Range("I45").Value = 1212.44
$ws.Range("J45").Value = 3625
$ws.Range("K45").Value = 1212.44
$ws.Range("L45").Value = 3625
$ws.Range("M45").Value = -835.4400000000001
$ws.Range("N45").Value = -4379
$ws.Range("H116").Value = 4903352.5
$ws.Range("I116").Value = 8929646
$ws.Range("J116").Value = 1776.9131
$ws.Range("K116").Value = 8929646
$ws.Range("L116").Value = 1776.9131
$ws.Range("M116").Value = -8927352
$ws.Range("N116").Value = -6364.9131
$ws.Range("H122").Value = 2279.5952
$ws.Range("I122").Value = 1886.3572
$ws.Range("J122").Value = 3066.0715
$ws.Range("K122").Value = 5659.071599999999
$ws.Range("L122").Value = 9198.2145
$ws.Range("M122").Value = -3209.071599999999
$ws.Range("N122").Value = -14098.2145
$ws.Range("H123").Value = 29995
$ws.Range("J123").Value = 29995
$ws.Range("L123").Value = 29995
$ws.Range("N123").Value = -39795
$ws.Range("H132").Value = 1993.5968
$ws.Range("I132").Value = 1636.0962
$ws.Range("J132").Value = 3852.6
$ws.Range("K132").Value = 4908.2886
$ws.Range("L132").Value = 11557.8
$ws.Range("M132").Value = -2378.2886
$ws.Range("N132").Value = -16617.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4903352.5
$ws.Range("I3").Value = 8929646
$ws.Range("J3").Value = 1776.9131
$ws.Range("K3").Value = 8929646
$ws.Range("L3").Value = 1776.9131
$ws.Range("M3").Value = -8929532
$ws.Range("N3").Value = -2004.9131
$ws.Range("H99").Value = 3858
$ws.Range("I99").Value = 3557.0908
$ws.Range("J99").Value = 4520
$ws.Range("K99").Value = 3557.0908
$ws.Range("L99").Value = 4520
$ws.Range("M99").Value = -2059.0908
$ws.Range("N99").Value = -7516
$ws.Range("H105").Value = 1365.762
$ws.Range("I105").Value = 1292.7778
$ws.Range("J105").Value = 1803.6666
$ws.Range("K105").Value = 1292.7778
$ws.Range("L105").Value = 1803.6666
$ws.Range("M105").Value = 454.2221999999999
$ws.Range("N105").Value = -5297.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1022988
$ws.Range("I31").Value = 1962475.6
$ws.Range("J31").Value = 3544
$ws.Range("K31").Value = 1962475.6
$ws.Range("L31").Value = 3544
$ws.Range("M31").Value = -1962180.6
$ws.Range("N31").Value = -4134
$ws.Range("H34").Value = 1022988
$ws.Range("I34").Value = 1962475.6
$ws.Range("J34").Value = 3544
$ws.Range("K34").Value = 1962475.6
$ws.Range("L34").Value = 3544
$ws.Range("M34").Value = -1962273.6
$ws.Range("N34").Value = -3948
$ws.Range("H105").Value = 2560.5833
$ws.Range("I105").Value = 2189.8572
$ws.Range("J105").Value = 3079.6
$ws.Range("K105").Value = 2189.8572
$ws.Range("L105").Value = 3079.6
$ws.Range("M105").Value = -442.8571999999999
$ws.Range("N105").Value = -6573.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3572300.2
$ws.Range("I113").Value = 25000516
$ws.Range("J113").Value = 931
$ws.Range("K113").Value = 75001548
$ws.Range("L113").Value = 2793
$ws.Range("M113").Value = -74999378
$ws.Range("N113").Value = -7133
$ws.Range("H118").Value = 2585.7144
$ws.Range("J118").Value = 2933.3333
$ws.Range("L118").Value = 8799.999899999999
$ws.Range("N118").Value = -11285.9999
$ws.Range("H131").Value = 1028.22
$ws.Range("I131").Value = 1347.9
$ws.Range("J131").Value = 992.7
$ws.Range("K131").Value = 4043.7
$ws.Range("L131").Value = 2978.1
$ws.Range("M131").Value = 996.2999999999997
$ws.Range("N131").Value = -13058.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5362.0938
$ws.Range("I122").Value = 4558.263
$ws.Range("J122").Value = 6536.923
$ws.Range("K122").Value = 13674.789
$ws.Range("L122").Value = 19610.769
$ws.Range("M122").Value = -11224.789
$ws.Range("N122").Value = -24510.769
$ws.Range("H133").Value = 28487.5
$ws.Range("J133").Value = 28487.5
$ws.Range("L133").Value = 28487.5
$ws.Range("N133").Value = -38607.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3495.25
$ws.Range("I68").Value = 1406.6666
$ws.Range("J68").Value = 5583.8335
$ws.Range("K68").Value = 1406.6666
$ws.Range("L68").Value = 5583.8335
$ws.Range("M68").Value = -657.6666
$ws.Range("N68").Value = -7081.8335
$ws.Range("H71").Value = 3495.25
$ws.Range("I71").Value = 1406.6666
$ws.Range("J71").Value = 5583.8335
$ws.Range("K71").Value = 7033.333000000001
$ws.Range("L71").Value = 27919.1675
$ws.Range("M71").Value = -3289.333000000001
$ws.Range("N71").Value = -35407.1675
$ws.Range("H122").Value = 2910.1428
$ws.Range("I122").Value = 2518.4614
$ws.Range("K122").Value = 7555.3842
$ws.Range("M122").Value = -5105.3842

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2384.7
$ws.Range("I122").Value = 1799.625
$ws.Range("K122").Value = 5398.875
$ws.Range("M122").Value = -2948.875
$ws.Range("I132").Value = 774281.4399999999
$ws.Range("J132").Value = 30483.777
$ws.Range("K132").Value = 2322844.32
$ws.Range("L132").Value = 91451.331
$ws.Range("M132").Value = -2320314.32
$ws.Range("N132").Value = -96511.331
